$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 65, shifting existing rows 65-105 down to 66-106
$ws.Rows.Item(65).Insert()

$ws.Cells.Item(65, 1).Value = 8
$ws.Cells.Item(65, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(65, 3).Value = "Coquimbo"
$ws.Cells.Item(65, 4).Value = 44438
$ws.Cells.Item(65, 5).Value = 4
$ws.Cells.Item(65, 6).Value = 100112037
$ws.Cells.Item(65, 7).Value = "Cebollín"
$ws.Cells.Item(65, 8).Value = "Sin especificar"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 3100
$ws.Cells.Item(65, 11).Value = 900
$ws.Cells.Item(65, 12).Value = 1000
$ws.Cells.Item(65, 13).Value = 950
$ws.Cells.Item(65, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(65, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(65, 16).Value = 158
$ws.Cells.Item(65, 17).Value = 6
$ws.Cells.Item(65, 18).Value = "Hortaliza"

$ws.Cells.Item(65, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
